$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: replace the broken "-----w" formula test case with a real test case ---
# New string order matters for the shared-strings table, so write C26, then B26, then D26.
$ws.Range("C26").Value = "The field being at the top should of the list of a contacts fields should only have 1 line above"
$ws.Range("B26").Value = "A contact does not have any associated numbers."
$ws.Range("D26").Value = "The top field has double lines and a longer distance to the start of the list than the number field would have."
$ws.Range("F26").Value = 0

# --- Fill in "Actual outcome" column (D) for the other JUnit test rows ---
$ws.Range("D10").Value = "Information shows up correctly."
$ws.Range("D9").Value = "Information copies successfully."
$ws.Range("D11").Value = "Information saves correctly"
$ws.Range("D12").Value = "Information saves correctly."
$ws.Range("D13").Value = "Information is deleted correctly."
$ws.Range("D14").Value = "Information is deleted correctly."
$ws.Range("D15").Value = "Application quits."
$ws.Range("D16").Value = "Application minimises."
$ws.Range("D17").Value = "Application maximises."
$ws.Range("D18").Value = "Windows minimum size is set and functions correctly."
$ws.Range("D19").Value = "Filepicker pops up and uploads pictures correctly"
$ws.Range("D20").Value = "Filepicker has a file filter for pictures."
$ws.Range("D21").Value = "Profile picture is deleted correctly."
$ws.Range("D22").Value = "User is prohibited of selected a different contact or leaving edit mode."
$ws.Range("D23").Value = "Field is added and focus is set."
$ws.Range("D24").Value = "Profile picture is deleted."
$ws.Range("D25").Value = "Contact is deleted."

# --- Cosmetic sheet-view changes ---
$ws.Application.ActiveWindow.Zoom = 107
$null = $ws.Range("B31").Select()

# --- Column F (hidden helper column) widened slightly ---
$ws.Columns(6).ColumnWidth = 10.498697916666666
